$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Map of 1-indexed table row -> new cell text (single-column table).
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "206"
    6  = "0.00030"
    7  = "0.00004"
    8  = "0.00002"
    9  = "0.00004"
    11 = "0.00008"
    12 = "0.00860"
    44 = "99.99"
    45 = "0.01"
    46 = "101"
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $tbl.Cell($rowIndex, 1)
    $cell.Range.Text = $updates[$rowIndex]
}
